$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: B1 1 -> 0, C1 0 -> 1
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1

# Row 2: A2 0 -> 1, C2 1 -> 0
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 0

# Row 3: A3 1 -> 0, B3 0 -> 1
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1
